$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 28575
$ws.Range("J3").Value = 28575
$ws.Range("L3").Value = 28575
$ws.Range("N3").Value = -28803
$ws.Range("H11").Value = 122.5
$ws.Range("I11").Value = 122.5
$ws.Range("K11").Value = 122.5
$ws.Range("M11").Value = 17.5
$ws.Range("H17").Value = 1726372.5
$ws.Range("J17").Value = 1756648.1
$ws.Range("L17").Value = 5269944.300000001
$ws.Range("N17").Value = -5270280.300000001
$ws.Range("H33").Value = 60.64706
$ws.Range("I33").Value = 62.4375
$ws.Range("J33").Value = 32
$ws.Range("K33").Value = 62.4375
$ws.Range("L33").Value = 32
$ws.Range("M33").Value = 166.5625
$ws.Range("N33").Value = -490
$ws.Range("H40").Value = 2256.125
$ws.Range("J40").Value = 2149.8572
$ws.Range("L40").Value = 2149.8572
$ws.Range("N40").Value = -2499.8572
$ws.Range("H102").Value = 28575
$ws.Range("J102").Value = 28575
$ws.Range("L102").Value = 28575
$ws.Range("N102").Value = -35065
$ws.Range("H129").Value = 204892.25
$ws.Range("J129").Value = 213594.05
$ws.Range("L129").Value = 640782.1499999999
$ws.Range("N129").Value = -650782.1499999999
$ws.Range("H132").Value = 2177.9788
$ws.Range("I132").Value = 2125.2144
$ws.Range("J132").Value = 2621.2
$ws.Range("K132").Value = 6375.6432
$ws.Range("L132").Value = 7863.599999999999
$ws.Range("M132").Value = -3845.6432
$ws.Range("N132").Value = -12923.6
$ws.Range("H135").Value = 20008174
$ws.Range("I135").Value = 1019.41174
$ws.Range("K135").Value = 9174.70566
$ws.Range("M135").Value = -6639.70566

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1393.7667
$ws.Range("I2").Value = 1325.05
$ws.Range("K2").Value = 1325.05
$ws.Range("M2").Value = -1212.05
$ws.Range("H32").Value = 24119.66
$ws.Range("I32").Value = 25479.303
$ws.Range("J32").Value = 9503.5
$ws.Range("K32").Value = 25479.303
$ws.Range("L32").Value = 9503.5
$ws.Range("M32").Value = -25192.303
$ws.Range("N32").Value = -10077.5
$ws.Range("H45").Value = 2765.2942
$ws.Range("I45").Value = 3412
$ws.Range("J45").Value = 2190.4443
$ws.Range("K45").Value = 3412
$ws.Range("L45").Value = 2190.4443
$ws.Range("M45").Value = -3035
$ws.Range("N45").Value = -2944.4443
$ws.Range("H61").Value = 2680.652
$ws.Range("I61").Value = 2097.5
$ws.Range("J61").Value = 4780
$ws.Range("K61").Value = 2097.5
$ws.Range("L61").Value = 4780
$ws.Range("M61").Value = -1885.5
$ws.Range("N61").Value = -5204
$ws.Range("H116").Value = 1393.7667
$ws.Range("I116").Value = 1325.05
$ws.Range("K116").Value = 1325.05
$ws.Range("M116").Value = 968.95
$ws.Range("H132").Value = 13080.467
$ws.Range("I132").Value = 1605.4849
$ws.Range("J132").Value = 44636.668
$ws.Range("K132").Value = 4816.4547
$ws.Range("L132").Value = 133910.004
$ws.Range("M132").Value = -2286.4547
$ws.Range("N132").Value = -138970.004
$ws.Range("H136").Value = 2680.652
$ws.Range("I136").Value = 2097.5
$ws.Range("J136").Value = 4780
$ws.Range("K136").Value = 6292.5
$ws.Range("L136").Value = 14340
$ws.Range("M136").Value = -3742.5
$ws.Range("N136").Value = -19440

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1393.7667
$ws.Range("I3").Value = 1325.05
$ws.Range("K3").Value = 1325.05
$ws.Range("M3").Value = -1211.05
$ws.Range("H33").Value = 10000
$ws.Range("J33").Value = 10000
$ws.Range("L33").Value = 10000
$ws.Range("N33").Value = -10672
$ws.Range("H99").Value = 2500
$ws.Range("I99").Value = 2333.3333
$ws.Range("J99").Value = 2625
$ws.Range("K99").Value = 2333.3333
$ws.Range("L99").Value = 2625
$ws.Range("M99").Value = -835.3332999999998
$ws.Range("N99").Value = -5621
$ws.Range("H100").Value = 22114.2
$ws.Range("J100").Value = 22114.2
$ws.Range("L100").Value = 22114.2
$ws.Range("N100").Value = -24278.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8801.046
$ws.Range("I31").Value = 10963.2
$ws.Range("J31").Value = 4167.857
$ws.Range("K31").Value = 10963.2
$ws.Range("L31").Value = 4167.857
$ws.Range("M31").Value = -10668.2
$ws.Range("N31").Value = -4757.857
$ws.Range("H34").Value = 8801.046
$ws.Range("I34").Value = 10963.2
$ws.Range("J34").Value = 4167.857
$ws.Range("K34").Value = 10963.2
$ws.Range("L34").Value = 4167.857
$ws.Range("M34").Value = -10761.2
$ws.Range("N34").Value = -4571.857
$ws.Range("H58").Value = 13835.205
$ws.Range("I58").Value = 1027.9286
$ws.Range("J58").Value = 46435.547
$ws.Range("K58").Value = 1027.9286
$ws.Range("L58").Value = 46435.547
$ws.Range("M58").Value = -824.9286
$ws.Range("N58").Value = -46841.547
$ws.Range("H74").Value = 33231.25
$ws.Range("J74").Value = 33231.25
$ws.Range("L74").Value = 33231.25
$ws.Range("N74").Value = -34979.25
$ws.Range("H77").Value = 33231.25
$ws.Range("J77").Value = 33231.25
$ws.Range("L77").Value = 99693.75
$ws.Range("N77").Value = -108429.75
$ws.Range("H92").Value = 29999.5
$ws.Range("J92").Value = 29999.5
$ws.Range("L92").Value = 29999.5
$ws.Range("N92").Value = -34991.5
$ws.Range("H96").Value = 11290.333
$ws.Range("J96").Value = 11290.333
$ws.Range("L96").Value = 11290.333
$ws.Range("N96").Value = -16782.333
$ws.Range("H132").Value = 20553.107
$ws.Range("I132").Value = 29402
$ws.Range("J132").Value = 4625.1
$ws.Range("K132").Value = 88206
$ws.Range("L132").Value = 13875.3
$ws.Range("M132").Value = -85676
$ws.Range("N132").Value = -18935.3
$ws.Range("H134").Value = 1368.9429
$ws.Range("I134").Value = 1013.26666
$ws.Range("J134").Value = 1635.7
$ws.Range("K134").Value = 3039.79998
$ws.Range("L134").Value = 4907.1
$ws.Range("M134").Value = -504.7999799999998
$ws.Range("N134").Value = -9977.1
$ws.Range("H136").Value = 13835.205
$ws.Range("I136").Value = 1027.9286
$ws.Range("J136").Value = 46435.547
$ws.Range("K136").Value = 3083.7858
$ws.Range("L136").Value = 139306.641
$ws.Range("M136").Value = -533.7857999999997
$ws.Range("N136").Value = -144406.641

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 1726.6666
$ws.Range("J112").Value = 3316.5
$ws.Range("L112").Value = 9949.5
$ws.Range("N112").Value = -12165.5
$ws.Range("H131").Value = 732.3099999999999
$ws.Range("J131").Value = 745.0213
$ws.Range("L131").Value = 2235.0639
$ws.Range("N131").Value = -12315.0639
$ws.Range("H132").Value = 1491.9
$ws.Range("I132").Value = 1369.8
$ws.Range("J132").Value = 1614
$ws.Range("K132").Value = 12328.2
$ws.Range("L132").Value = 14526
$ws.Range("M132").Value = -9798.199999999999
$ws.Range("N132").Value = -19586

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 18442
$ws.Range("J45").Value = 18442
$ws.Range("L45").Value = 18442
$ws.Range("N45").Value = -19560
$ws.Range("H132").Value = 106516.6
$ws.Range("I132").Value = 119415.22
$ws.Range("J132").Value = 87168.664
$ws.Range("K132").Value = 358245.66
$ws.Range("L132").Value = 261505.992
$ws.Range("M132").Value = -355715.66
$ws.Range("N132").Value = -266565.992

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1216.6666
$ws.Range("I46").Value = 775
$ws.Range("J46").Value = 2100
$ws.Range("K46").Value = 775
$ws.Range("L46").Value = 2100
$ws.Range("M46").Value = -587
$ws.Range("N46").Value = -2476
$ws.Range("H82").Value = 3115.25
$ws.Range("I82").Value = 3257.1428
$ws.Range("J82").Value = 2916.6
$ws.Range("K82").Value = 3257.1428
$ws.Range("L82").Value = 2916.6
$ws.Range("M82").Value = -2896.1428
$ws.Range("N82").Value = -3638.6
$ws.Range("H85").Value = 3115.25
$ws.Range("I85").Value = 3257.1428
$ws.Range("J85").Value = 2916.6
$ws.Range("K85").Value = 3257.1428
$ws.Range("L85").Value = 2916.6
$ws.Range("M85").Value = -2009.1428
$ws.Range("N85").Value = -5412.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1726.3334
$ws.Range("I122").Value = 1723.6086
$ws.Range("J122").Value = 1735.2858
$ws.Range("K122").Value = 5170.825800000001
$ws.Range("L122").Value = 5205.857400000001
$ws.Range("M122").Value = -2720.825800000001
$ws.Range("N122").Value = -10105.8574
